$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the computed values in column J (recomputed "base opt" plot data) ---
$ws.Range("J1").Value = 15.109654900000001
$ws.Range("J2").Value = 104.45106560000001
$ws.Range("J3").Value = 113.1285275
$ws.Range("J4").Value = 200.41288249999999

# --- Re-apply column widths (the widths were recalculated/tightened in the refreshed export) ---
# ColumnWidth is expressed in characters; the host snaps it to a whole-pixel grid
# (pixels = round(width * 6); stored = (pixels + 5) / 6), so we back-solve for the
# character width that reproduces each target stored width as closely as the grid allows.
$ws.Columns.Item(1).ColumnWidth  = 2.1666666666666665   # -> 3
$ws.Columns.Item(2).ColumnWidth  = 4.166666666666667    # -> 5
$ws.Columns.Item(3).ColumnWidth  = 1.1666666666666667   # -> 2
$ws.Columns.Item(4).ColumnWidth  = 2.6666666666666665   # -> 3.5546875 (closest: 3.5)
$ws.Columns.Item(5).ColumnWidth  = 2.1666666666666665   # -> 3
$ws.Columns.Item(6).ColumnWidth  = 3.1666666666666665   # -> 4
$ws.Columns.Item(7).ColumnWidth  = 3.6666666666666665   # -> 4.5546875 (closest: 4.5)
$ws.Columns.Item(8).ColumnWidth  = 2.1666666666666665   # -> 3
$ws.Columns.Item(9).ColumnWidth  = 3.1666666666666665   # -> 4
$ws.Columns.Item(10).ColumnWidth = 10.666666666666666   # -> 11.5546875 (closest: 11.5)
$ws.Columns.Item(11).ColumnWidth = 14.666666666666666   # -> 15.5546875 (closest: 15.5)
$ws.Columns.Item(12).ColumnWidth = 14.333333333333334   # -> 15.21875 (closest: 15.1666...)
$ws.Columns.Item(13).ColumnWidth = 1.1666666666666667   # -> 2
$ws.Columns.Item(14).ColumnWidth = 2.1666666666666665   # -> 3
